$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Header row rename
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# 2. Title-case Spanish connector words in municipality/state names
$textChanges = @(
    @('B4', 'Pabellón De Arteaga'),
    @('B5', 'Rincón De Romos'),
    @('B6', 'San José De Gracia'),
    @('B25', 'Amatenango De La Frontera'),
    @('B29', 'Comitán De Domínguez'),
    @('B46', 'Ocozocoautla De Espinosa'),
    @('B77', 'Guadalupe Y Calvo'),
    @('B78', 'Hidalgo Del Parral'),
    @('B86', 'San Francisco Del Oro'),
    @('B88', 'Valle De Zaragoza'),
    @('B104', 'Villa De Álvarez'),
    @('A106', 'Ciudad De México'),
    @('B136', 'San Juan Del Río'),
    @('A144', 'Estado De México'),
    @('B144', 'Acambay De Ruíz Castañeda'),
    @('B147', 'Almoloya De Juárez'),
    @('B148', 'Almoloya Del Río'),
    @('B152', 'Atizapán De Zaragoza'),
    @('B157', 'Chapa De Mota'),
    @('B159', 'Coacalco De Berriozábal'),
    @('B164', 'Ecatepec De Morelos'),
    @('B168', 'Ixtapan De La Sal'),
    @('B169', 'Ixtapan Del Oro'),
    @('B178', 'Naucalpan De Juárez'),
    @('B185', 'San Felipe Del Progreso'),
    @('B186', 'San Martín De Las Pirámides'),
    @('B204', 'Tlalnepantla De Baz'),
    @('B210', 'Valle De Chalco Solidaridad'),
    @('B211', 'Villa De Allende'),
    @('B212', 'Villa Del Carbón'),
    @('B223', 'San Miguel De Allende'),
    @('B224', 'Apaseo El Alto'),
    @('B225', 'Apaseo El Grande'),
    @('B232', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B236', 'Jaral Del Progreso'),
    @('B242', 'Purísima Del Rincón'),
    @('B247', 'San Francisco Del Rincón'),
    @('B249', 'San Luis De La Paz'),
    @('B250', 'Santa Cruz De Juventino Rosas'),
    @('B251', 'Silao De La Victoria'),
    @('B255', 'Valle De Santiago'),
    @('B260', 'Acapulco De Juárez'),
    @('B262', 'Ajuchitlán Del Progreso'),
    @('B263', 'Alcozauca De Guerrero'),
    @('B268', 'Atoyac De Álvarez'),
    @('B269', 'Ayutla De Los Libres'),
    @('B272', 'Buenavista De Cuéllar'),
    @('B273', 'Chilapa De Álvarez'),
    @('B274', 'Chilpancingo De Los Bravo'),
    @('B275', 'Coahuayutla De José María Izazaga'),
    @('B280', 'Coyuca De Benítez'),
    @('B281', 'Coyuca De Catalán'),
    @('B285', 'Cuetzala Del Progreso'),
    @('B286', 'Cutzamala De Pinzón'),
    @('B292', 'Huitzuco De Los Figueroa'),
    @('B293', 'Iguala De La Independencia'),
    @('B295', 'Ixcateopan De Cuauhtémoc'),
    @('B296', 'Zihuatanejo De Azueta'),
    @('B298', 'La Unión De Isidoro Montes De Oca'),
    @('B301', 'Mártir De Cuilapan'),
    @('B313', 'Taxco De Alarcón'),
    @('B315', 'Técpan De Galeana'),
    @('B317', 'Tepecoacuilco De Trujano'),
    @('B318', 'Tixtla De Guerrero'),
    @('B321', 'Tlalixtaquilla De Maldonado'),
    @('B322', 'Tlapa De Comonfort'),
    @('B332', 'Atotonilco El Grande'),
    @('B338', 'Cuautepec De Hinojosa'),
    @('B341', 'Huasca De Ocampo'),
    @('B343', 'Huejutla De Reyes'),
    @('B350', 'Mineral Del Chico'),
    @('B351', 'Mixquiahuala De Juárez'),
    @('B352', 'Molango De Escamilla'),
    @('B353', 'Nopala De Villagrán'),
    @('B354', 'Omitlán De Juárez'),
    @('B355', 'Pachuca De Soto'),
    @('B356', 'Progreso De Obregón'),
    @('B358', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B361', 'Tenango De Doria'),
    @('B363', 'Tepehuacán De Guerrero'),
    @('B364', 'Tepeji Del Río De Ocampo'),
    @('B365', 'Tezontepec De Aldama'),
    @('B371', 'Tula De Allende'),
    @('B372', 'Tulancingo De Bravo'),
    @('B376', 'Acatlán De Juárez'),
    @('B377', 'Ahualulco De Mercado'),
    @('B381', 'Atemajac De Brizuela'),
    @('B384', 'Atotonilco El Alto'),
    @('B386', 'Autlán De Navarro'),
    @('B391', 'Cañadas De Obregón'),
    @('B398', 'Concepción De Buenos Aires'),
    @('B399', 'Cuautitlán De García Barragán'),
    @('B406', 'Encarnación De Díaz'),
    @('B412', 'Huejuquilla El Alto'),
    @('B413', 'Ixtlahuacán De Los Membrillos'),
    @('B414', 'Ixtlahuacán Del Río'),
    @('B418', 'Jilotlán De Los Dolores'),
    @('B424', 'Lagos De Moreno'),
    @('B431', 'Ojuelos De Jalisco'),
    @('B436', 'San Cristóbal De La Barranca'),
    @('B438', 'San Juan De Los Lagos'),
    @('B441', 'San Martín De Bolaños'),
    @('B443', 'San Miguel El Alto'),
    @('B444', 'San Sebastián Del Oeste'),
    @('B447', 'Talpa De Allende'),
    @('B448', 'Tamazula De Gordiano'),
    @('B453', 'Teocuitatlán De Corona'),
    @('B454', 'Tepatitlán De Morelos'),
    @('B456', 'Tizapán El Alto'),
    @('B457', 'Tlajomulco De Zúñiga'),
    @('B469', 'Unión De San Antonio'),
    @('B470', 'Unión De Tula'),
    @('B471', 'Valle De Guadalupe'),
    @('B475', 'Yahualica De González Gallo'),
    @('B476', 'Zacoalco De Torres'),
    @('B479', 'Zapotitlán De Vadillo'),
    @('B480', 'Zapotlán Del Rey'),
    @('B481', 'Zapotlán El Grande'),
    @('B504', 'Coalcomán De Vázquez Pallares'),
    @('B506', 'Cojumatlán De Régules'),
    @('B568', 'Tiquicheo De Nicolás Romero'),
    @('B598', 'Jonacatepec De Leandro Valle'),
    @('B602', 'Puente De Ixtla'),
    @('B607', 'Tetela Del Volcán'),
    @('B608', 'Tlaltizapán De Zapata'),
    @('B618', 'Amatlán De Cañas'),
    @('B619', 'Bahía De Banderas'),
    @('B623', 'Ixtlán Del Río'),
    @('B630', 'Santa María Del Oro'),
    @('B642', 'Mier Y Noriega'),
    @('B645', 'Acatlán De Pérez Figueroa'),
    @('B651', 'Chalcatongo De Hidalgo'),
    @('B654', 'Coicoyán De Las Flores'),
    @('B656', 'Constancia Del Rosario'),
    @('B658', 'Cuilápam De Guerrero'),
    @('B659', 'Guadalupe De Ramírez'),
    @('B660', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B661', 'Heroica Ciudad De Huajuapan De León'),
    @('B662', 'Heroica Ciudad De Tlaxiaco'),
    @('B664', 'Ixtlán De Juárez'),
    @('B665', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B670', 'Mariscala De Juárez'),
    @('B673', 'Miahuatlán De Porfirio Díaz'),
    @('B674', 'Nejapa De Madero'),
    @('B675', 'Oaxaca De Juárez'),
    @('B676', 'Ocotlán De Morelos'),
    @('B677', 'Putla Villa De Guerrero'),
    @('B682', 'San Agustín De Las Juntas'),
    @('B689', 'San Antonio De La Cal'),
    @('B715', 'San Juan Del Río'),
    @('B742', 'San Miguel El Grande'),
    @('B751', 'San Pedro El Alto'),
    @('B758', 'San Pedro Y San Pablo Teposcolula'),
    @('B765', 'Santa Ana Del Valle'),
    @('B773', 'Santa Cruz Tacache De Mina'),
    @('B777', 'Santa Inés Del Monte'),
    @('B778', 'Santa Lucía Del Camino'),
    @('B797', 'Santiago Del Río'),
    @('B815', 'Santo Domingo De Morelos'),
    @('B822', 'Tamazulápam Del Espíritu Santo'),
    @('B823', 'Teotitlán De Flores Magón'),
    @('B824', 'Teotitlán Del Valle'),
    @('B825', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B826', 'Tlacolula De Matamoros'),
    @('B827', 'Villa De Chilapa De Díaz'),
    @('B828', 'Villa De Etla'),
    @('B829', 'Villa De Tututepec'),
    @('B830', 'Villa De Zaachila'),
    @('B832', 'Villa Sola De Vega'),
    @('B834', 'Zapotitlán Del Río'),
    @('B836', 'Zimatlán De Álvarez'),
    @('B861', 'Cuetzalan Del Progreso'),
    @('B874', 'Izúcar De Matamoros'),
    @('B879', 'Mazapiltepec De Juárez'),
    @('B884', 'Palmar De Bravo'),
    @('B897', 'San Salvador El Seco'),
    @('B898', 'San Salvador El Verde'),
    @('B901', 'Tecali De Herrera'),
    @('B908', 'Tepanco De López'),
    @('B909', 'Tepatlaxco De Hidalgo'),
    @('B912', 'Tepexi De Rodríguez'),
    @('B914', 'Tetela De Ocampo'),
    @('B918', 'Tlacotepec De Benito Juárez'),
    @('B927', 'Totoltepec De Guerrero'),
    @('B944', 'Amealco De Bonfil'),
    @('B946', 'Cadereyta De Montes'),
    @('B951', 'Landa De Matamoros'),
    @('B953', 'Pinal De Amoles'),
    @('B956', 'San Juan Del Río'),
    @('B971', 'Mexquitic De Carmona'),
    @('B973', 'San Ciro De Acosta'),
    @('B975', 'Santa María Del Río'),
    @('B981', 'Villa De Arista'),
    @('B982', 'Villa De Ramos'),
    @('B983', 'Villa De Reyes'),
    @('B1036', 'Soto La Marina'),
    @('B1046', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1048', 'San Pablo Del Monte'),
    @('B1050', 'Tepetitla De Lardizábal'),
    @('B1070', 'Amatlán De Los Reyes'),
    @('B1082', 'Cazones De Herrera'),
    @('B1090', 'Cosamaloapan De Carpio'),
    @('B1103', 'Hueyapan De Ocampo'),
    @('B1105', 'Ixhuatlán Del Café'),
    @('B1112', 'Juchique De Ferrer'),
    @('B1116', 'Lerdo De Tejada'),
    @('B1120', 'Martínez De La Torre'),
    @('B1121', 'Medellín De Bravo'),
    @('B1124', 'Mixtla De Altamirano'),
    @('B1132', 'Paso De Ovejas'),
    @('B1134', 'Poza Rica De Hidalgo'),
    @('B1140', 'Sayula De Alemán'),
    @('B1144', 'Soledad De Doblado'),
    @('B1165', 'Vega De Alatorre'),
    @('B1170', 'Zozocolco De Hidalgo'),
    @('B1200', 'Concepción Del Oro'),
    @('B1202', 'El Plateado De Joaquín Amaro'),
    @('B1216', 'Moyahua De Estrada'),
    @('B1217', 'Nochistlán De Mejía'),
    @('B1225', 'Teúl De González Ortega'),
    @('B1226', 'Tlaltenango De Sánchez Román'),
    @('B1227', 'Trinidad García De La Cadena'),
    @('B1229', 'Villa De Cos')
)
foreach ($pair in $textChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# 3. Floating point precision fix in column D
$dCells = @(
    'D95',
    'D104',
    'D178',
    'D185',
    'D227',
    'D267',
    'D280',
    'D324',
    'D332',
    'D349',
    'D426',
    'D493',
    'D532',
    'D543',
    'D591',
    'D623',
    'D645',
    'D653',
    'D662',
    'D664',
    'D671',
    'D858',
    'D902',
    'D938',
    'D949',
    'D997',
    'D1088'
)
foreach ($cell in $dCells) {
    $ws.Range($cell).Value = 0.0009415309292910272
}

# 4. Remove trailing footer/metadata rows (1237:1241)
$ws.Rows('1237:1241').Delete()